$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.951.94"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.46%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.641.95"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.13%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.10%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.08%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5043"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.75%  "

# Row 7
$ws.Range("E7").Value = "  -0.19%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2572"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.44%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06414"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.03%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.60"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.24%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07778"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.71%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.267"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.08%  "

# Row 13
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.869.74"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.01%  "

# Row 14
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.553.31"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.58%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5439"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.42%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅7932"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.62%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.43"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.00%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.986.44"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.31%  "

# Row 19
$ws.Range("E19").Value = "  -0.24%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "199.64"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.59%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.380"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.41%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.901"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.10%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.976"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.10%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.004"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.20%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.891"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.04%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "140.67"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.50%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1135"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.38%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.819"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.87%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.66"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.06%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.242"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.21%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.04921"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.23%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.263"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.25%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.210"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.26%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.540"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.23%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.373"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.50%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.8924"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.19%  "

# Row 37
$ws.Range("B37").Value = "MXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.606"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.70%  "

# Row 38
$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.154.06"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.92%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5562"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.37%  "

# Row 40
$ws.Range("E40").Value = "  +0.18%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.005"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.08%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.718"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.61%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8107"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.11%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.74"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.04%  "

# Row 45
$ws.Range("E45").Value = "  +4.81%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.780.72"
$ws.Range("D46").Style = "Normal"

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4521"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.07%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.005"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.16%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "54.79"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.48%  "

# Row 50
$ws.Range("E50").Value = "  -0.01%  "

# Row 51
$ws.Range("E51").Value = "  -0.04%  "
